# edit.ps1
# Applies the commit "add new test cases of the login" changes:
#  1. Turn the blank paragraph right after the title into a (still empty)
#     paragraph whose paragraph mark carries Bold/size-32 run formatting,
#     and insert a brand-new paragraph after it with the "Linked In :"
#     text (including the proofErr gramStart/gramEnd markers Word's
#     grammar checker adds around "In :").
#  2. Move the <w:lastRenderedPageBreak/> marker from the "Suspension and
#     resumption criteria" bullet to the "Pricing" bullet.

$d = $word.ActiveDocument

$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Change 1: blank paragraph after the title -----------------------------
# Find the empty paragraph immediately following the "Manual Testing
# Project" heading paragraph.
$headingParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "Manual Testing Project") {
        $headingParaIndex = $i
        break
    }
}
$blankPara = $d.Paragraphs.Item($headingParaIndex + 1)

$blankParaXml = $pkgHeader + '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr></w:p>' + $pkgFooter
[void]$blankPara.Range.InsertXML($blankParaXml)

# Re-fetch: the (still blank) paragraph is now Paragraphs.Item($headingParaIndex + 1)
$blankPara = $d.Paragraphs.Item($headingParaIndex + 1)
$blankPara.Range.InsertParagraphAfter()

$linkedInPara = $d.Paragraphs.Item($headingParaIndex + 2)
$linkedInXml = $pkgHeader + '<w:p><w:r><w:t xml:space="preserve">Linked </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>In :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>www.linkedin.com/in/shubham-randive</w:t></w:r></w:p>' + $pkgFooter
[void]$linkedInPara.Range.InsertXML($linkedInXml)

# --- Change 2: move lastRenderedPageBreak from "Suspension..." to "Pricing" -
$pricingIndex = 0
$suspensionIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "^Pricing") {
        $pricingIndex = $i
    }
    if ($t -match "^Suspension and resumption criteria") {
        $suspensionIndex = $i
    }
}

$pPricing = $d.Paragraphs.Item($pricingIndex)
$pricingXml = $pkgHeader + '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Pricing</w:t></w:r></w:p>' + $pkgFooter
[void]$pPricing.Range.InsertXML($pricingXml)

$pSuspension = $d.Paragraphs.Item($suspensionIndex)
$suspensionXml = $pkgHeader + '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Suspension and resumption criteria</w:t></w:r></w:p>' + $pkgFooter
[void]$pSuspension.Range.InsertXML($suspensionXml)
